$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "Software"
$ws.Range("J3").Value = "Software"

$ws.Range("J3").Select()
